$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row into the translation table at sheet row 32 ---
# (Table "Tabla13" currently spans B2:F203; a new row is inserted as the
#  31st data row, shifting all subsequent rows down by one.)
$ws.Range("B32:F32").Insert()

# Populate the newly inserted row with the new "strWindowPos" resource entry
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E32").Value = "Remember window position and size on startup"
$ws.Range("F32").Value = ""

# Existing row 25 (strChkDlgPath) now also carries the same UI-location comment
$ws.Range("D25").Value = 'In "settings" form, tab "User interface"'

# Grow the table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))

# Widen column D slightly to fit the new/longer comment text
$ws.Columns.Item(4).ColumnWidth = 34.8
